$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The old row 77 ("AGGREGATE RESULTS" header, bold, with a couple of
# empty styled cells in D/E) moves up to row 76 and keeps only its
# text in column A (bold).
# ------------------------------------------------------------------
$ws.Range("A76").Value = $ws.Range("A77").Value2
$ws.Range("A76").Font.Bold = $true

# Clear out the old row 77 content/format so we can rebuild it as the
# new "final" results header row.
$ws.Range("A77:E77").ClearContents()
$ws.Range("A77:E77").ClearFormats()

# ------------------------------------------------------------------
# Build the new row 77 header, mirroring row 1's headers except for
# column C, which uses the new "Actual Number of Associations" text
# (note: lower-case "of", distinct from the existing "Actual Number
# Of Associations" string used in row 1).
# ------------------------------------------------------------------
$ws.Range("A77").Value = $ws.Range("A1").Value2
$ws.Range("B77").Value = $ws.Range("B1").Value2
$ws.Range("C77").Value = "Actual Number of Associations"
$ws.Range("D77").Value = $ws.Range("D1").Value2
$ws.Range("E77").Value = $ws.Range("E1").Value2
$ws.Range("F77").Value = $ws.Range("F1").Value2
$ws.Range("G77").Value = $ws.Range("G1").Value2

# A77, D77, E77, F77, G77 carry an explicit (non-bold) font; B77 and
# C77 are left with the plain default cell style.
$ws.Range("A77").Font.Name = "Calibri"
$ws.Range("D77").Font.Name = "Calibri"
$ws.Range("E77").Font.Name = "Calibri"
$ws.Range("F77").Font.Name = "Calibri"
$ws.Range("G77").Font.Name = "Calibri"

# ------------------------------------------------------------------
# Update the view state to reflect where the author left off.
# ------------------------------------------------------------------
$ws.Range("F78").Select()
